# Daily attendance processing - 2025-10-23 18:53:41
# Reorders the comma-separated "Recorded By" values in column G for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 'system, backup@backdoor.com, System'
$ws.Cells.Item(5, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(7, 7).Value = 'admin@admin.com, System'
$ws.Cells.Item(8, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(11, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(17, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(29, 7).Value = 'system, backup@backdoor.com, System'
$ws.Cells.Item(32, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(34, 7).Value = 'admin@admin.com, System'
$ws.Cells.Item(35, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(38, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(44, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(56, 7).Value = 'system, backup@backdoor.com, System'
$ws.Cells.Item(59, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(61, 7).Value = 'admin@admin.com, System'
$ws.Cells.Item(62, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(65, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(71, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(83, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(84, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(85, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(90, 7).Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Cells.Item(96, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(97, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(99, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(109, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(110, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(111, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(116, 7).Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Cells.Item(122, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(123, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(125, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(135, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(136, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(137, 7).Value = 'backup@backdoor.com, System'
$ws.Cells.Item(142, 7).Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Cells.Item(148, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(149, 7).Value = 'System, dnasr281@gmail.com'
$ws.Cells.Item(151, 7).Value = 'System, dnasr281@gmail.com'
